$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete entire row 7 ("log_voting_power"), shifting all rows below it up by one.
$ws.Rows.Item(7).Delete()

# Update name and value in the row that is now row 3 (Rolling_Avg_Misaligned_1M -> _6M, H3 0.27 -> 0.25)
$ws.Range("A3").Value = "Rolling_Avg_Misaligned_6M"
$ws.Range("H3").Value = 0.25
